$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "2024-10-05"
$ws.Range("A5").NumberFormat = "General"
$ws.Range("B5").Value = 0.9959
